# Rename the three inline picture drawings (the two Pearson Edexcel logos
# in the primary/first-page footers and the BTEC logo in the first-page
# header) the same way Word itself exposes a rename: InlineShape has no
# writable .Name, so convert to a floating Shape, set .Name (which is
# serialized back out as the drawing's wp:docPr/@name), then convert back
# to an inline shape so the wp:inline wrapper (and everything else about
# the drawing) is left exactly as it was.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlinePicture($inlineShape, $newName) {
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape()
}

# Footer (default / primary) -> footer1.xml: Pearson logo image1.png -> image2.png
$footer1 = $sec.Footers.Item(1)
$shapes1 = $footer1.Range.InlineShapes
Rename-InlinePicture $shapes1.Item(1) "image2.png"

# Footer (first page) -> footer2.xml: Pearson logo image1.png -> image2.png
$footer2 = $sec.Footers.Item(2)
$shapes2 = $footer2.Range.InlineShapes
Rename-InlinePicture $shapes2.Item(1) "image2.png"

# Header (first page) -> header2.xml: BTEC logo image2.jpg -> image1.jpg
$header2 = $sec.Headers.Item(2)
$shapes3 = $header2.Range.InlineShapes
Rename-InlinePicture $shapes3.Item(1) "image1.jpg"
